# Fruta / hortaliza, semanal
# A new weekly price record was added to the "Coliflor" sheet. It belongs
# chronologically right after the existing row 195 (same market/category
# block), so we insert a brand-new row at position 195 - shifting every
# following row down by one (195-249 -> 196-250) - and then populate the
# newly inserted row with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 195; everything from old row 195 downward
# (through the former last row 249) shifts down to 196..250.
$ws.Rows.Item(195).Insert()

# Populate the newly inserted row 195 with the new record.
$ws.Cells.Item(195, 1).Value = 4
$ws.Cells.Item(195, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(195, 3).Value = "Los Lagos"
$ws.Cells.Item(195, 4).Value = 44551
$ws.Cells.Item(195, 5).Value = 10
$ws.Cells.Item(195, 6).Value = 100112008
$ws.Cells.Item(195, 7).Value = "Coliflor"
$ws.Cells.Item(195, 8).Value = "Sin especificar"
$ws.Cells.Item(195, 9).Value = "Primera"
$ws.Cells.Item(195, 10).Value = 1000
$ws.Cells.Item(195, 11).Value = 1200
$ws.Cells.Item(195, 12).Value = 1200
$ws.Cells.Item(195, 13).Value = 1200
$ws.Cells.Item(195, 14).Value = "`$/unidad"
$ws.Cells.Item(195, 15).Value = "Región Metropolitana"
$ws.Cells.Item(195, 16).Value = 1200
$ws.Cells.Item(195, 17).Value = 1
$ws.Cells.Item(195, 18).Value = "Hortaliza"
